$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row order/content (rows 2-10), reflecting the rearranged Viterbi results table.
# Columns: A=label, B..I = numeric probability values
$rows = @(
    @{ Row = 2;  Label = "NCMS000"; Values = @(0, 0, 0, 0.01017241781572992, 0, 0, 0, 0) },
    @{ Row = 3;  Label = "VMIP3S0"; Values = @(0.0012941074971961, 0, 0, 0, 0, 0, 0, 0) },
    @{ Row = 4;  Label = "NCMP000"; Values = @(0, 0, 0, 0, 0, 0, 0.0007549414348462665, 0) },
    @{ Row = 5;  Label = "AQ0MS0";  Values = @(0, 0, 0, 0.0009827570803180559, 0, 0, 0, 0) },
    @{ Row = 6;  Label = "NCFS000"; Values = @(0, 0, 0, 0, 0, 0, 0, 0) },
    @{ Row = 7;  Label = "SPS00";   Values = @(0, 0.01959015197765447, 0, 0, 0, 0.1249534018034177, 0, 0) },
    @{ Row = 8;  Label = "Fp";      Values = @(0, 0, 0, 0, 0, 0, 0, 0.0959409594095941) },
    @{ Row = 9;  Label = "AQ0CS0";  Values = @(0, 0, 0, 0, 0.002729616337259263, 0, 0, 0) },
    @{ Row = 10; Label = "DA0MS0";  Values = @(0, 0, 0.1398289673695107, 0, 0, 0, 0, 0) }
)

foreach ($entry in $rows) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.Label
    $colIndex = 2
    foreach ($val in $entry.Values) {
        $ws.Cells.Item($r, $colIndex).Value = $val
        $colIndex++
    }
}
